$d = $word.ActiveDocument

function ReplaceText($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# Returns the 1-based Paragraphs index whose trimmed text matches $text exactly.
function FindParagraphIndex($text) {
    $n = $d.Paragraphs.Count
    for ($i = 1; $i -le $n; $i++) {
        $t = $d.Paragraphs.Item($i).Range.Text
        $t = $t.TrimEnd([char]13, [char]7)
        if ($t -eq $text) {
            return $i
        }
    }
    return -1
}

# 1. Update generated timestamp
ReplaceText "Generated: 2026-02-09 14:13 " "Generated: 2026-02-09 14:44"

# 2. Overview paragraph: "two UIs: a user camera interface" -> "a user camera interface"
ReplaceText "This application implements a modular anomaly inspection pipeline with a web UI. It supports category-specific anomaly detection on MVTec AD-style data and offers two UIs: a user camera interface and an admin upload interface." "This application implements a modular anomaly inspection pipeline with a web UI. It supports category-specific anomaly detection on MVTec AD-style data and offers a user camera interface and an admin upload interface."

# 3. configs/base.yaml bullet
ReplaceText "- configs/base.yaml: Primary configuration file (categories, thresholds, labels, RPM, paths)." "- configs/base.yaml: Primary configuration file."

# 4-11. Project structure sub-bullets
ReplaceText "  - src/api.py: FastAPI server with user/admin UIs and upload/analyze endpoints." "  - src/api.py: FastAPI server with user/admin UIs."
ReplaceText "  - src/pipeline.py: Batch pipeline (baseline anomaly detection and post-processing)." "  - src/pipeline.py: Batch pipeline."
ReplaceText "  - src/models/mean_diff.py: Baseline anomaly model (mean-difference)." "  - src/models/mean_diff.py: Baseline anomaly model."
ReplaceText "  - src/data/mvtec.py: MVTec data loader/iterator." "  - src/data/mvtec.py: MVTec data loader."
ReplaceText "  - src/postproc/: Heatmap normalization, thresholding, bounding boxes." "  - src/postproc/: Heatmap, mask, and bbox utilities."
ReplaceText "  - src/vlm/semantics.py: Placeholder VLM integration (returns Unknown)." "  - src/vlm/semantics.py: LLaVA integration."
ReplaceText "  - src/risk/: RPM lookup + policy mapping placeholders." "  - src/risk/: RPM lookup + policy mapping."
ReplaceText "  - src/uncertainty/: confidence fusion and review rules." "  - src/uncertainty/: confidence and review rules."

# 12. data folder bullet
ReplaceText "- data/: Dataset folder (e.g., data/bottle)." "- data/: Dataset folder."

# 13. Heading "3. Configuration (configs/base.yaml)" -> "3. Configuration"
ReplaceText "3. Configuration (configs/base.yaml)" "3. Configuration"

# 14. "Important sections:" -> "Important sections in configs/base.yaml:"
ReplaceText "Important sections:" "Important sections in configs/base.yaml:"

# 15-18. Config section bullets
ReplaceText "- postproc: heatmap normalization, thresholding, min area, image-level threshold." "- postproc: heatmap normalization and thresholding."
ReplaceText "- labels: category-specific label sets (fixed)." "- labels: category-specific label sets."
ReplaceText "- risk: RPM table and risk-to-action mapping (placeholders until fully defined)." "- risk: RPM table and risk-to-action mapping."
ReplaceText "- uncertainty: confidence fusion method and review threshold." "- uncertainty: confidence fusion and review threshold."

# 19. FastAPI exposes line
ReplaceText "The FastAPI app exposes two UIs and API endpoints:" "The FastAPI app exposes:"

# 20-21. Upload/Analyze endpoint bullets
ReplaceText "- Upload endpoint: POST /upload (saves an uploaded image)." "- Upload endpoint: POST /upload."
ReplaceText "- Analyze endpoint: POST /analyze (runs baseline anomaly detection)." "- Analyze endpoint: POST /analyze."

# 22. Baseline description paragraph
ReplaceText "The current baseline uses a mean-difference model: it computes the per-pixel mean and std from training good images and scores anomalies based on normalized absolute deviation." "Current baseline uses a mean-difference model: per-pixel mean and std are computed from training images and anomalies are scored by normalized deviation."

# 23. Delete the "Outputs:" paragraph and the 4 bullet paragraphs that follow it
$outputsIdx = FindParagraphIndex "Outputs:"
$outputsStart = $d.Paragraphs.Item($outputsIdx).Range.Start
$outputsEnd = $d.Paragraphs.Item($outputsIdx + 4).Range.End
$d.Range($outputsStart, $outputsEnd).Delete()

# 24. VLM/Risk/Uncertainty paragraph
ReplaceText "VLM semantics, RPM lookup, and action policy are currently stubbed. The system returns Unknown labels and no risk until the RPM and label sets are fully populated, and VLM integration is implemented." "VLM semantics are implemented with LLaVA-1.6 (Mistral). Risk uses a deterministic RPM lookup, and uncertainty combines anomaly and VLM confidence."

# 25. Delete "8. Notes and Limitations" and "9. Suggested Next Steps" sections (everything
# after the "python -m uvicorn src.api:app --reload" paragraph)
$lastKeptIdx = FindParagraphIndex "python -m uvicorn src.api:app --reload"
$n = $d.Paragraphs.Count
$tailStart = $d.Paragraphs.Item($lastKeptIdx + 1).Range.Start
$tailEnd = $d.Paragraphs.Item($n).Range.End
$d.Range($tailStart, $tailEnd).Delete()
